$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# --- Numeric value updates ---
$ws.Range("C16").Value = 1
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 63
$ws.Range("K16").Value = 85.294117647058
$ws.Range("L16").Value = 53.658536585365
$ws.Range("M16").Value = -35.051546391752
$ws.Range("N16").Value = -83.507853403141
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 98
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = 8.888888888888
$ws.Range("L17").Value = 22.5
$ws.Range("M17").Value = 12.643678160919
$ws.Range("N17").Value = -66.552901023890
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 30
$ws.Range("J18").Value = 72
$ws.Range("K18").Value = 51.388888888888
$ws.Range("M18").Value = 10.101010101010
$ws.Range("N18").Value = -73.543689320388
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 8
$ws.Range("H19").Value = 187.5
$ws.Range("I19").Value = 208
$ws.Range("J19").Value = 121
$ws.Range("K19").Value = 71.900826446281
$ws.Range("L19").Value = 41.496598639455
$ws.Range("M19").Value = -4.587155963302
$ws.Range("N19").Value = -7.555555555555
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = -75
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = -5.263157894736
$ws.Range("N20").Value = -86.010362694300
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = -61.538461538461
$ws.Range("F21").Value = 59
$ws.Range("G21").Value = 39
$ws.Range("H21").Value = 51.282051282051
$ws.Range("I21").Value = 540
$ws.Range("J21").Value = 381
$ws.Range("K21").Value = 41.732283464566
$ws.Range("L21").Value = 30.750605326876
$ws.Range("M21").Value = -3.914590747330
$ws.Range("N21").Value = -68.568102444703
$ws.Range("M22").Value = -46.153846153846
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 16.666666666666
$ws.Range("I23").Value = 111
$ws.Range("J23").Value = 88
$ws.Range("K23").Value = 26.136363636363
$ws.Range("L23").Value = 33.734939759036
$ws.Range("M23").Value = 48
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 40
$ws.Range("G24").Value = 26
$ws.Range("H24").Value = 53.846153846153
$ws.Range("I24").Value = 437
$ws.Range("J24").Value = 336
$ws.Range("K24").Value = 30.059523809523
$ws.Range("L24").Value = 4.545454545454
$ws.Range("M24").Value = 5.048076923076
$ws.Range("C25").Value = 1
$ws.Range("E25").Value = -75
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 162
$ws.Range("J25").Value = 151
$ws.Range("K25").Value = 7.284768211920
$ws.Range("L25").Value = 47.272727272727
$ws.Range("M25").Value = -40.659340659340
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("L27").Value = 144.444444444444
$ws.Range("D28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = -60
$ws.Range("L28").Value = -60
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -69.230769230769
$ws.Range("L29").Value = -71.428571428571

# --- Cells that become "no data" text markers ("0" / "***.*") ---
# Copy style+value from existing donor cells that already hold these shared strings
# so the exact same shared-string index & cell style are reused.
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
